$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the branch name text (A3): Floorspace -> Dwellings
$ws.Range("A3").Value = "CIMS.CAN.QC.Residential.Dwellings.Lighting"

# Row 3 values M3:W3 - replace formulas with static values
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 0.5
$ws.Range("Q3").Value = 0.0001
$ws.Range("R3").Value = 0.0001
$ws.Range("S3").Value = 0.0001
$ws.Range("T3").Value = 0.0001
$ws.Range("U3").Value = 0.0001
$ws.Range("V3").Value = 0.0001
$ws.Range("W3").Value = 0.0001

# Update the selection to A1:X4
[void]$ws.Range("A1:X4").Select()
